$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 3 and row 4 text/number values (A:F) first ---
# Row 3 (mirrors row 2 contents for A-F)
$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = "Admin"
$ws.Range("C3").Value = "IND_DAU_51"
$ws.Range("D3").Value = "IDM+18"
$ws.Range("E3").Value = "10.75.58.51"
$ws.Range("F3").Value = 409026540

# Row 4 (mirrors row 2 contents for A-F)
$ws.Range("A4").Value = "Admin"
$ws.Range("B4").Value = "Admin"
$ws.Range("C4").Value = "IND_DAU_51"
$ws.Range("D4").Value = "IDM+18"
$ws.Range("E4").Value = "10.75.58.51"
$ws.Range("F4").Value = 409026540

# --- Update the numeric-looking text values (stored as text via shared strings) ---
# Order matters so the resulting shared-string table matches the target layout:
# new unique strings should appear in order: 500, 1200, 400, 10000, 1000 (30000 reused)
$ws.Range("K2").Value = "'500"
$ws.Range("K4").Value = "'1200"
$ws.Range("G2").Value = "'400"
$ws.Range("G3").Value = "'10000"
$ws.Range("K3").Value = "'30000"
$ws.Range("G4").Value = "'1000"

$ws.Range("H3").Value = "'530"
$ws.Range("I3").Value = "'60"
$ws.Range("J3").Value = "'1"

$ws.Range("H4").Value = "'530"
$ws.Range("I4").Value = "'60"
$ws.Range("J4").Value = "'1"

# --- Apply the same cell formatting (quote-prefix style) used in row 2 to rows 3 & 4 ---
$ws.Range("F2:K2").Copy()
$ws.Range("F3:K3").PasteSpecial(-4122)
$ws.Range("F2:K2").Copy()
$ws.Range("F4:K4").PasteSpecial(-4122)

# --- Update selection to match target view state ---
$ws.Range("K4").Select()
